$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 501, shifting existing rows 501:525 down to 502:526
$ws.Rows(501).Insert()

# Populate the new row 501 with data (mirrors the surrounding rows, with updated values)
$ws.Range("A501").Value = 9
$ws.Range("B501").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C501").Value = "Metropolitana"
$ws.Range("D501").Value = 45267
$ws.Range("E501").Value = 13
$ws.Range("F501").Value = 300000001
$ws.Range("G501").Value = "Rabanito"
$ws.Range("H501").Value = "Sin especificar"
$ws.Range("I501").Value = "Primera"
$ws.Range("J501").Value = 7000
$ws.Range("K501").Value = 3000
$ws.Range("L501").Value = 3000
$ws.Range("M501").Value = 3000
$ws.Range("N501").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O501").Value = "Provincia de Chacabuco"
$ws.Range("P501").Value = 30
$ws.Range("Q501").Value = 100
$ws.Range("R501").Value = "Hortaliza"
